$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the registered-family row (row 4) figures for 2015-2021 ---
$ws.Range("E4").Value = 7005
$ws.Range("F4").Value = 4674
$ws.Range("G4").Value = 4333
$ws.Range("H4").Value = 4191
$ws.Range("I4").Value = 3987
$ws.Range("J4").Value = 4166
$ws.Range("K4").Value = 4376

# --- Update the subsistence-allowance row (row 5) figures for 2015-2021 ---
$ws.Range("E5").Value = 2040
$ws.Range("F5").Value = 1984
$ws.Range("G5").Value = 1712
$ws.Range("H5").Value = 1572
$ws.Range("I5").Value = 1423
$ws.Range("J5").Value = 1682
$ws.Range("K5").Value = 2047

# --- Residual cell-style artifact left behind by the source data refresh
#     (pasting the corrected figures in from another workbook pulls that
#     workbook's default "Normal" style in as "Normal 2") ---
$ghostStyle = $wb.Styles.Add("Normal 2")
$ghostStyle.Font.Name = "Arial"
